$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 298.92856
$ws.Range("I19").Value = 187.77777
$ws.Range("J19").Value = 351.57895
$ws.Range("K19").Value = 187.77777
$ws.Range("L19").Value = 351.57895
$ws.Range("M19").Value = -12.77777
$ws.Range("N19").Value = -701.5789500000001

# Row 32
$ws.Range("H32").Value = 8131923.5
$ws.Range("I32").Value = 30304418
$ws.Range("J32").Value = 2009.0667
$ws.Range("K32").Value = 30304418
$ws.Range("L32").Value = 2009.0667
$ws.Range("M32").Value = -30304092
$ws.Range("N32").Value = -2661.0667

# Row 38
$ws.Range("H38").Value = 935.7
$ws.Range("I38").Value = 509.5
$ws.Range("J38").Value = 1575
$ws.Range("K38").Value = 1528.5
$ws.Range("L38").Value = 4725
$ws.Range("M38").Value = -1156.5
$ws.Range("N38").Value = -5469

# Row 46
$ws.Range("H46").Value = 1441.625
$ws.Range("I46").Value = 1802.8
$ws.Range("K46").Value = 5408.4
$ws.Range("M46").Value = -5289.4

# Row 60
$ws.Range("H60").Value = 1441.625
$ws.Range("I60").Value = 1802.8
$ws.Range("K60").Value = 5408.4
$ws.Range("M60").Value = -4924.4

# Row 112
$ws.Range("H112").Value = 5001024
$ws.Range("J112").Value = 6250955
$ws.Range("L112").Value = 18752865
$ws.Range("N112").Value = -18755081

# Row 129
$ws.Range("H129").Value = 1054.7046
$ws.Range("I129").Value = 2998.25
$ws.Range("J129").Value = 860.35
$ws.Range("K129").Value = 8994.75
$ws.Range("L129").Value = 2581.05
$ws.Range("M129").Value = -3994.75
$ws.Range("N129").Value = -12581.05

# Row 139
$ws.Range("H139").Value = 48499.6
$ws.Range("J139").Value = 48499.6
$ws.Range("L139").Value = 48499.6
$ws.Range("N139").Value = -58779.6


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2884.923
$ws.Range("I61").Value = 2334.3333
$ws.Range("J61").Value = 4123.75
$ws.Range("K61").Value = 2334.3333
$ws.Range("L61").Value = 4123.75
$ws.Range("M61").Value = -2122.3333
$ws.Range("N61").Value = -4547.75

# Row 110
$ws.Range("H110").Value = 4480.5
$ws.Range("I110").Value = 1887.4166
$ws.Range("K110").Value = 1887.4166
$ws.Range("M110").Value = 157.5834

# Row 122
$ws.Range("H122").Value = 1973.6842
$ws.Range("I122").Value = 1980.6666
$ws.Range("J122").Value = 1947.5
$ws.Range("K122").Value = 5941.9998
$ws.Range("L122").Value = 5842.5
$ws.Range("M122").Value = -3491.9998
$ws.Range("N122").Value = -10742.5

# Row 136
$ws.Range("H136").Value = 2884.923
$ws.Range("I136").Value = 2334.3333
$ws.Range("J136").Value = 4123.75
$ws.Range("K136").Value = 7002.999899999999
$ws.Range("L136").Value = 12371.25
$ws.Range("M136").Value = -4452.999899999999
$ws.Range("N136").Value = -17471.25


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 94
$ws.Range("H94").Value = 1005.56525
$ws.Range("I94").Value = 378.66666
$ws.Range("K94").Value = 378.66666
$ws.Range("M94").Value = 72.33334000000002

# Row 99
$ws.Range("H99").Value = 2991.818
$ws.Range("I99").Value = 3056.6667
$ws.Range("J99").Value = 2700
$ws.Range("K99").Value = 3056.6667
$ws.Range("L99").Value = 2700
$ws.Range("M99").Value = -1558.6667
$ws.Range("N99").Value = -5696

# Row 122
$ws.Range("H122").Value = 300737.5
$ws.Range("I122").Value = 300737.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 902212.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -899762.5
$ws.Range("N122").ClearContents()

# Row 126
$ws.Range("H126").Value = 2991.818
$ws.Range("I126").Value = 3056.6667
$ws.Range("J126").Value = 2700
$ws.Range("K126").Value = 9170.000100000001
$ws.Range("L126").Value = 8100
$ws.Range("M126").Value = -6700.000100000001
$ws.Range("N126").Value = -13040


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5116.364
$ws.Range("I5").Value = 7441.4287
$ws.Range("J5").Value = 1047.5
$ws.Range("K5").Value = 22324.2861
$ws.Range("L5").Value = 3142.5
$ws.Range("M5").Value = -22212.2861
$ws.Range("N5").Value = -3366.5

# Row 132
$ws.Range("H132").Value = 2013.5264
$ws.Range("I132").Value = 1067.25
$ws.Range("J132").Value = 3635.7144
$ws.Range("K132").Value = 9605.25
$ws.Range("L132").Value = 32721.4296
$ws.Range("M132").Value = -7075.25
$ws.Range("N132").Value = -37781.4296

# Row 133
$ws.Range("H133").Value = 4182.5
$ws.Range("I133").Value = 4182.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 12547.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -7487.5
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 59011060
$ws.Range("I134").Value = 125390600
$ws.Range("K134").Value = 376171800
$ws.Range("M134").Value = -376166730

# Row 135
$ws.Range("H135").Value = 5116.364
$ws.Range("I135").Value = 7441.4287
$ws.Range("J135").Value = 1047.5
$ws.Range("K135").Value = 66972.85830000001
$ws.Range("L135").Value = 9427.5
$ws.Range("M135").Value = -64437.85830000001
$ws.Range("N135").Value = -14497.5

# Row 136
$ws.Range("H136").Value = 55557404
$ws.Range("I136").Value = 83334720
$ws.Range("J136").Value = 2777.6667
$ws.Range("K136").Value = 250004160
$ws.Range("L136").Value = 8333.000100000001
$ws.Range("M136").Value = -249999060
$ws.Range("N136").Value = -18533.0001

# Row 137
$ws.Range("H137").Value = 8798.786
$ws.Range("J137").Value = 11916.625
$ws.Range("L137").Value = 35749.875
$ws.Range("N137").Value = -45949.875

# Row 138
$ws.Range("H138").Value = 1314.75
$ws.Range("I138").Value = 1314.75
$ws.Range("K138").Value = 3944.25
$ws.Range("M138").Value = 1195.75

# Row 139
$ws.Range("H139").Value = 86100.47
$ws.Range("I139").Value = 189556.19
$ws.Range("J139").Value = 3335.9
$ws.Range("K139").Value = 568668.5700000001
$ws.Range("L139").Value = 10007.7
$ws.Range("M139").Value = -563528.5700000001
$ws.Range("N139").Value = -20287.7

# Row 140
$ws.Range("H140").Value = 2830
$ws.Range("I140").Value = 859.6
$ws.Range("J140").Value = 4472
$ws.Range("K140").Value = 2578.8
$ws.Range("L140").Value = 13416
$ws.Range("M140").Value = 2601.2
$ws.Range("N140").Value = -23776

# Row 141
$ws.Range("H141").Value = 125004770
$ws.Range("I141").Value = 250003620
$ws.Range("J141").Value = 5925
$ws.Range("K141").Value = 750010860
$ws.Range("L141").Value = 17775
$ws.Range("M141").Value = -750005680
$ws.Range("N141").Value = -28135


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3975
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 4633.3335
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 4633.3335
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -7877.3335

# Row 122
$ws.Range("H122").Value = 2866.1333
$ws.Range("I122").Value = 2382.7693
$ws.Range("J122").Value = 6008
$ws.Range("K122").Value = 7148.3079
$ws.Range("L122").Value = 18024
$ws.Range("M122").Value = -4698.3079
$ws.Range("N122").Value = -22924

# Row 126
$ws.Range("H126").Value = 14836.479
$ws.Range("I126").Value = 28830.818
$ws.Range("J126").Value = 2008.3334
$ws.Range("K126").Value = 86492.454
$ws.Range("L126").Value = 6025.0002
$ws.Range("M126").Value = -84022.454
$ws.Range("N126").Value = -10965.0002


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1266.5555
$ws.Range("I22").Value = 1119.8
$ws.Range("K22").Value = 1119.8
$ws.Range("M22").Value = -824.8

# Row 27
$ws.Range("H27").Value = 1266.5555
$ws.Range("I27").Value = 1119.8
$ws.Range("K27").Value = 1119.8
$ws.Range("M27").Value = -1012.8

# Row 40
$ws.Range("H40").Value = 4368.737
$ws.Range("I40").Value = 3646.4614
$ws.Range("J40").Value = 5933.6665
$ws.Range("K40").Value = 3646.4614
$ws.Range("L40").Value = 5933.6665
$ws.Range("M40").Value = -3510.4614
$ws.Range("N40").Value = -6205.6665

# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 100
$ws.Range("H100").Value = 2746.4614
$ws.Range("I100").Value = 1531.4286
$ws.Range("J100").Value = 4164
$ws.Range("K100").Value = 1531.4286
$ws.Range("L100").Value = 4164
$ws.Range("M100").Value = -990.4286
$ws.Range("N100").Value = -5246


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 43000
$ws.Range("J93").Value = 43000
$ws.Range("L93").Value = 43000
$ws.Range("N93").Value = -47992

# Row 122
$ws.Range("H122").Value = 1036.3
$ws.Range("I122").Value = 1036.3
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3108.9
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -658.8999999999996
$ws.Range("N122").ClearContents()

